$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1066.3625
$ws.Range("J17").Value = 1066.3948
$ws.Range("L17").Value = 3199.1844
$ws.Range("N17").Value = -3535.1844
$ws.Range("H132").Value = 2027.5227
$ws.Range("I132").Value = 1128.3334
$ws.Range("J132").Value = 5293
$ws.Range("K132").Value = 3385.0002
$ws.Range("L132").Value = 15879
$ws.Range("M132").Value = -855.0001999999999
$ws.Range("N132").Value = -20939
$ws.Range("H135").Value = 486.6087
$ws.Range("I135").Value = 383.42856
$ws.Range("J135").Value = 1570
$ws.Range("K135").Value = 3450.85704
$ws.Range("L135").Value = 14130
$ws.Range("M135").Value = -915.8570399999999
$ws.Range("N135").Value = -19200
$ws.Range("H137").Value = 2664.4
$ws.Range("I137").Value = 2423.4075
$ws.Range("J137").Value = 4833.3335
$ws.Range("K137").Value = 7270.2225
$ws.Range("L137").Value = 14500.0005
$ws.Range("M137").Value = -4720.2225
$ws.Range("N137").Value = -19600.0005
$ws.Range("H141").Value = 812.0961
$ws.Range("I141").Value = 663.08887
$ws.Range("J141").Value = 1770
$ws.Range("K141").Value = 1989.26661
$ws.Range("L141").Value = 5310
$ws.Range("M141").Value = 3190.73339
$ws.Range("N141").Value = -15670

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6051.0806
$ws.Range("I32").Value = 3822.1428
$ws.Range("J32").Value = 31405.25
$ws.Range("K32").Value = 3822.1428
$ws.Range("L32").Value = 31405.25
$ws.Range("M32").Value = -3535.1428
$ws.Range("N32").Value = -31979.25
$ws.Range("H74").Value = 152302.4
$ws.Range("I74").Value = 193115.83
$ws.Range("J74").Value = 51240.617
$ws.Range("K74").Value = 193115.83
$ws.Range("L74").Value = 51240.617
$ws.Range("M74").Value = -192241.83
$ws.Range("N74").Value = -52988.617
$ws.Range("H77").Value = 152302.4
$ws.Range("I77").Value = 193115.83
$ws.Range("J77").Value = 51240.617
$ws.Range("K77").Value = 965579.1499999999
$ws.Range("L77").Value = 256203.085
$ws.Range("M77").Value = -961211.1499999999
$ws.Range("N77").Value = -264939.085
$ws.Range("H102").Value = 1897.5358
$ws.Range("I102").Value = 1582.5
$ws.Range("J102").Value = 2317.5833
$ws.Range("K102").Value = 1582.5
$ws.Range("L102").Value = 2317.5833
$ws.Range("M102").Value = 39.5
$ws.Range("N102").Value = -5561.5833
$ws.Range("H122").Value = 6320.7856
$ws.Range("I122").Value = 6299.1
$ws.Range("K122").Value = 18897.3
$ws.Range("M122").Value = -16447.3
$ws.Range("H132").Value = 1131.54
$ws.Range("I132").Value = 842.0278
$ws.Range("J132").Value = 1876
$ws.Range("K132").Value = 2526.0834
$ws.Range("L132").Value = 5628
$ws.Range("M132").Value = 3.916600000000017
$ws.Range("N132").Value = -10688

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1752.6
$ws.Range("I20").Value = 1683
$ws.Range("J20").Value = 1807.2858
$ws.Range("K20").Value = 1683
$ws.Range("L20").Value = 1807.2858
$ws.Range("M20").Value = -1436
$ws.Range("N20").Value = -2301.2858
$ws.Range("H86").Value = 3681.0908
$ws.Range("I86").Value = 4712.2666
$ws.Range("J86").Value = 1471.4286
$ws.Range("K86").Value = 4712.2666
$ws.Range("L86").Value = 1471.4286
$ws.Range("M86").Value = -3589.2666
$ws.Range("N86").Value = -3717.4286
$ws.Range("H89").Value = 3681.0908
$ws.Range("I89").Value = 4712.2666
$ws.Range("J89").Value = 1471.4286
$ws.Range("K89").Value = 23561.333
$ws.Range("L89").Value = 7357.143
$ws.Range("M89").Value = -17945.333
$ws.Range("N89").Value = -18589.143
$ws.Range("H105").Value = 2030.0526
$ws.Range("I105").Value = 1904.2858
$ws.Range("J105").Value = 2382.2
$ws.Range("K105").Value = 1904.2858
$ws.Range("L105").Value = 2382.2
$ws.Range("M105").Value = -157.2858000000001
$ws.Range("N105").Value = -5876.2
$ws.Range("H134").Value = 2195.9858
$ws.Range("I134").Value = 1967.6938
$ws.Range("J134").Value = 2704.4546
$ws.Range("K134").Value = 5903.0814
$ws.Range("L134").Value = 8113.3638
$ws.Range("M134").Value = -3368.0814
$ws.Range("N134").Value = -13183.3638

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2983.4888
$ws.Range("I31").Value = 1714.3871
$ws.Range("J31").Value = 5793.643
$ws.Range("K31").Value = 1714.3871
$ws.Range("L31").Value = 5793.643
$ws.Range("M31").Value = -1419.3871
$ws.Range("N31").Value = -6383.643
$ws.Range("H34").Value = 2983.4888
$ws.Range("I34").Value = 1714.3871
$ws.Range("J34").Value = 5793.643
$ws.Range("K34").Value = 1714.3871
$ws.Range("L34").Value = 5793.643
$ws.Range("M34").Value = -1512.3871
$ws.Range("N34").Value = -6197.643
$ws.Range("H132").Value = 1207.1625
$ws.Range("I132").Value = 699.5574
$ws.Range("J132").Value = 2836.842
$ws.Range("K132").Value = 2098.6722
$ws.Range("L132").Value = 8510.526
$ws.Range("M132").Value = 431.3278
$ws.Range("N132").Value = -13570.526

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 2096.6667
$ws.Range("I133").Value = 2096.6667
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 6290.000100000001
$ws.Range("L133").Value = 0
$ws.Range("M133").Value = -1230.000100000001
$ws.Range("N133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2599.9443
$ws.Range("I132").Value = 2060.319
$ws.Range("J132").Value = 3614.44
$ws.Range("K132").Value = 6180.957
$ws.Range("L132").Value = 10843.32
$ws.Range("M132").Value = -3650.957
$ws.Range("N132").Value = -15903.32

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2671.5
$ws.Range("I7").Value = 2682.5
$ws.Range("K7").Value = 2682.5
$ws.Range("M7").Value = -2570.5
$ws.Range("H100").Value = 71433880
$ws.Range("I100").Value = 11800
$ws.Range("J100").Value = 111112810
$ws.Range("K100").Value = 11800
$ws.Range("L100").Value = 111112810
$ws.Range("M100").Value = -11259
$ws.Range("N100").Value = -111113892
$ws.Range("H122").Value = 1984.9375
$ws.Range("I122").Value = 1772
$ws.Range("J122").Value = 2623.75
$ws.Range("K122").Value = 5316
$ws.Range("L122").Value = 7871.25
$ws.Range("M122").Value = -2866
$ws.Range("N122").Value = -12771.25
$ws.Range("H126").Value = 2671.5
$ws.Range("I126").Value = 2682.5
$ws.Range("K126").Value = 8047.5
$ws.Range("M126").Value = -5577.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1267.5
$ws.Range("I100").Value = 2225.8
$ws.Range("J100").Value = 831.9091
$ws.Range("K100").Value = 4451.6
$ws.Range("L100").Value = 1663.8182
$ws.Range("M100").Value = -3910.6
$ws.Range("N100").Value = -2745.8182
$ws.Range("H105").Value = 36003.625
$ws.Range("J105").Value = 36003.625
$ws.Range("L105").Value = 36003.625
$ws.Range("N105").Value = -42991.625
$ws.Range("H132").Value = 1749.4606
$ws.Range("I132").Value = 1001.6739
$ws.Range("J132").Value = 2896.0667
$ws.Range("K132").Value = 3005.0217
$ws.Range("L132").Value = 8688.2001
$ws.Range("M132").Value = -475.0217000000002
$ws.Range("N132").Value = -13748.2001
